$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2719.5
$ws.Range("I18").Value = 2719.5
$ws.Range("K18").Value = 2719.5
$ws.Range("M18").Value = -2435.5

$ws.Range("H49").Value = 2000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 2000
$ws.Range("K49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("M49").Value = 6000
$ws.Range("N49").Value = -6272

$ws.Range("H51").Value = 9750
$ws.Range("I51").Value = 9750
$ws.Range("K51").Value = 9750
$ws.Range("M51").Value = -9266

$ws.Range("H58").Value = 406.14285
$ws.Range("I58").Value = 12.666667
$ws.Range("J58").Value = 701.25
$ws.Range("K58").Value = 38.000001
$ws.Range("L58").Value = 2103.75
$ws.Range("M58").Value = 111.999999
$ws.Range("N58").Value = -2403.75

$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14126

$ws.Range("H70").Value = 9183.333000000001
$ws.Range("I70").Value = 15002
$ws.Range("K70").Value = 45006
$ws.Range("M70").Value = -44736

$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40632

$ws.Range("H73").Value = 9183.333000000001
$ws.Range("I73").Value = 15002
$ws.Range("K73").Value = 45006
$ws.Range("M73").Value = -44070

$ws.Range("H113").Value = 2800
$ws.Range("I113").Value = 2752.5
$ws.Range("J113").Value = 2895
$ws.Range("K113").Value = 2752.5
$ws.Range("L113").Value = 2895
$ws.Range("M113").Value = 501.5
$ws.Range("N113").Value = -9403

$ws.Range("H115").Value = 975
$ws.Range("J115").Value = 450
$ws.Range("L115").Value = 1350
$ws.Range("N115").Value = -4484

$ws.Range("H116").Value = 4624.6665
$ws.Range("J116").Value = 4624.6665
$ws.Range("L116").Value = 4624.6665
$ws.Range("N116").Value = -11508.6665

$ws.Range("H118").Value = 189.5
$ws.Range("I118").Value = 189.5
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 568.5
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = 1088.5

$ws.Range("H129").Value = 9419
$ws.Range("I129").Value = 9419
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 28257
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -23257

$ws.Range("H138").Value = 2426.5
$ws.Range("J138").Value = 3729.1667
$ws.Range("L138").Value = 11187.5001
$ws.Range("N138").Value = -21467.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2582.8333
$ws.Range("I88").Value = 2166
$ws.Range("J88").Value = 2999.6667
$ws.Range("K88").Value = 2166
$ws.Range("L88").Value = 2999.6667
$ws.Range("M88").Value = -1760
$ws.Range("N88").Value = -3811.6667

$ws.Range("H91").Value = 2582.8333
$ws.Range("I91").Value = 2166
$ws.Range("J91").Value = 2999.6667
$ws.Range("K91").Value = 2166
$ws.Range("L91").Value = 2999.6667
$ws.Range("M91").Value = -762
$ws.Range("N91").Value = -5807.6667

$ws.Range("H98").Value = 54982.168
$ws.Range("J98").Value = 54982.168
$ws.Range("L98").Value = 54982.168
$ws.Range("N98").Value = -60972.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("N15").Value = 0

$ws.Range("H82").Value = 4625
$ws.Range("I82").Value = 4625
$ws.Range("K82").Value = 4625
$ws.Range("M82").Value = -4242

$ws.Range("H85").Value = 4625
$ws.Range("I85").Value = 4625
$ws.Range("K85").Value = 4625
$ws.Range("M85").Value = -3299

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 273.8
$ws.Range("I2").Value = 193.5
$ws.Range("K2").Value = 193.5
$ws.Range("M2").Value = -80.5

$ws.Range("H7").Value = 209.72917
$ws.Range("I7").Value = 207.6923
$ws.Range("J7").Value = 210.48572
$ws.Range("K7").Value = 207.6923
$ws.Range("L7").Value = 210.48572
$ws.Range("M7").Value = -94.69229999999999
$ws.Range("N7").Value = -436.48572

$ws.Range("H17").Value = 208
$ws.Range("I17").Value = 208
$ws.Range("K17").Value = 208
$ws.Range("M17").Value = -34

$ws.Range("H22").Value = 102
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H55").Value = 11795
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H59").Value = 17497.5
$ws.Range("I59").Value = 17497.5
$ws.Range("K59").Value = 17497.5
$ws.Range("M59").Value = -16352.5

$ws.Range("H86").Value = 6698.8
$ws.Range("I86").Value = 5248.75
$ws.Range("K86").Value = 5248.75
$ws.Range("M86").Value = -4125.75

$ws.Range("H89").Value = 6698.8
$ws.Range("I89").Value = 5248.75
$ws.Range("K89").Value = 26243.75
$ws.Range("M89").Value = -20627.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 142864860
$ws.Range("I4").Value = 819.4
$ws.Range("K4").Value = 2458.2
$ws.Range("M4").Value = -2346.2

$ws.Range("H17").Value = 42.333332
$ws.Range("J17").Value = 42
$ws.Range("L17").Value = 126
$ws.Range("N17").Value = -464

$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("N19").Value = 0

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H118").Value = 3264.5
$ws.Range("I118").Value = 3264.5
$ws.Range("K118").Value = 9793.5
$ws.Range("M118").Value = -8550.5

$ws.Range("H119").Value = 1000
$ws.Range("I119").Value = 1000
$ws.Range("K119").Value = 3000
$ws.Range("M119").Value = 1838

$ws.Range("H129").Value = 2632.5
$ws.Range("J129").Value = 2632.5
$ws.Range("L129").Value = 7897.5
$ws.Range("N129").Value = -17897.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 39.6
$ws.Range("I2").Value = 24.625
$ws.Range("J2").Value = 99.5
$ws.Range("K2").Value = 24.625
$ws.Range("L2").Value = 99.5
$ws.Range("M2").Value = 88.375
$ws.Range("N2").Value = -325.5

$ws.Range("H80").Value = 7899.2
$ws.Range("J80").Value = 8749
$ws.Range("L80").Value = 8749
$ws.Range("N80").Value = -10745

$ws.Range("H83").Value = 7899.2
$ws.Range("J83").Value = 8749
$ws.Range("L83").Value = 43745
$ws.Range("N83").Value = -53729

$ws.Range("H102").Value = 1917.9166
$ws.Range("I102").Value = 1917.9166
$ws.Range("K102").Value = 1917.9166
$ws.Range("M102").Value = -295.9166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 2468
$ws.Range("I10").Value = 1750
$ws.Range("J10").Value = 3904
$ws.Range("K10").Value = 1750
$ws.Range("L10").Value = 3904
$ws.Range("M10").Value = -1610
$ws.Range("N10").Value = -4184

$ws.Range("H22").Value = 3453.2222
$ws.Range("J22").Value = 3874.25
$ws.Range("L22").Value = 3874.25
$ws.Range("N22").Value = -4464.25

$ws.Range("H27").Value = 3453.2222
$ws.Range("J27").Value = 3874.25
$ws.Range("L27").Value = 3874.25
$ws.Range("N27").Value = -4088.25

$ws.Range("H40").Value = 7455.5
$ws.Range("I40").Value = 3676.25
$ws.Range("K40").Value = 3676.25
$ws.Range("M40").Value = -3540.25

$ws.Range("H55").Value = 2585
$ws.Range("I55").Value = 2421.7144
$ws.Range("K55").Value = 2421.7144
$ws.Range("M55").Value = -2248.7144

$ws.Range("H132").Value = 9595
$ws.Range("I132").Value = 9514.556
$ws.Range("J132").Value = 9836.333000000001
$ws.Range("K132").Value = 28543.668
$ws.Range("L132").Value = 29508.999
$ws.Range("M132").Value = -26013.668
$ws.Range("N132").Value = -34568.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H62").Value = 10749.5
$ws.Range("I62").Value = 10749.5
$ws.Range("K62").Value = 10749.5
$ws.Range("M62").Value = -10125.5

$ws.Range("H65").Value = 10749.5
$ws.Range("I65").Value = 10749.5
$ws.Range("K65").Value = 53747.5
$ws.Range("M65").Value = -50627.5

$ws.Range("H107").Value = 1700.8572
$ws.Range("I107").Value = 1580
$ws.Range("K107").Value = 4740
$ws.Range("M107").Value = -2820
